$p = $ppt.ActivePresentation

# --- Slide 1 (SlideID 256): reassign the shape currently using Id=1 ---
# PowerPoint shape Ids are read-only via automation; the only supported way
# to make the engine hand out a fresh sequential Id (1 -> 4, since the
# current max Id on the slide is 3) is to copy the shape, delete the
# original, and paste the copy back -- then restore its original z-order
# position so the shape tree order is unchanged.
$s1 = $p.Slides.Item(1)

$target = $null
for ($j = 1; $j -le $s1.Shapes.Count; $j++) {
    $cand = $s1.Shapes.Item($j)
    if ($cand.Id -eq 1) {
        $target = $cand
        break
    }
}

if ($target -ne $null) {
    $origPos = $target.ZOrderPosition
    $target.Copy()
    $target.Delete()
    $newShape = $s1.Shapes.Paste()
    # The pasted shape lands at the end of the z-order; walk it back to the
    # position the original shape occupied.
    $stepsBack = $newShape.ZOrderPosition - $origPos
    for ($k = 0; $k -lt $stepsBack; $k++) {
        $newShape.ZOrder(3)   # msoSendBackward
    }
}

# --- Remove slide 2 (SlideID 284) entirely ---
$p.Slides.Item(2).Delete()
